$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 157, shifting the existing weekly records
# (formerly rows 157-187) down to rows 158-188.
$ws.Rows("157:157").Insert()

# Populate the newly inserted row 157 with the new weekly record.
$ws.Cells.Item(157, 1).Value  = 10
$ws.Cells.Item(157, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(157, 3).Value  = "La Araucanía"
$ws.Cells.Item(157, 4).Value  = 45015
$ws.Cells.Item(157, 5).Value  = 9
$ws.Cells.Item(157, 6).Value  = 100114002
$ws.Cells.Item(157, 7).Value  = "Camote"
$ws.Cells.Item(157, 8).Value  = "Sin especificar"
$ws.Cells.Item(157, 9).Value  = "Primera"
$ws.Cells.Item(157, 10).Value = 60
$ws.Cells.Item(157, 11).Value = 26000
$ws.Cells.Item(157, 12).Value = 26000
$ws.Cells.Item(157, 13).Value = 26000
$ws.Cells.Item(157, 14).Value = "$/malla 20 kilos"
$ws.Cells.Item(157, 15).Value = "Perú"
$ws.Cells.Item(157, 16).Value = 1300
$ws.Cells.Item(157, 17).Value = 20
$ws.Cells.Item(157, 18).Value = "Hortaliza"

# Make sure the date cell keeps the same date/time number format used by
# the rest of column D.
$ws.Cells.Item(157, 4).NumberFormat = $ws.Cells.Item(158, 4).NumberFormat
